$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (Firm, Time, Lawyers Registered) replacing the old rows 2-26
# with 23 new rows (rows 2-24). Row count shrinks from 26 to 24.
$data = @(
  @("Havel Partners", "17s", "1"),
  @("Blakes", "51s", "1"),
  @("Mijares Angoitia Cortés And Fuentes", "07s", "1"),
  @("Hill Dickinson", "08s", "2"),
  @("McCarthy Tetrault", "09s", "1"),
  @("Bennett Jones", "18s", "1"),
  @("Howse Williams", "09s", "1"),
  @("Jones Day", "21s", "3"),
  @("Cassels", "01min 05s", "1"),
  @("Ritch Mueller And Nicolau", "06s", "1"),
  @("Krogerus", "07s", "1"),
  @("Latham And Watkins", "20s", "3"),
  @("Kinstellar", "42s", "12"),
  @("Lee And Ko", "19s", "1"),
  @("Stikeman Elliott", "01min 19s", "0"),
  @("Hannes Snellman", "04s", "1"),
  @("Peter And Kim", "05s", "2"),
  @("Kromann Reumert", "06s", "1"),
  @("Hakun Law", "05s", "1"),
  @("Dittmar And Indrenius", "22s", "0"),
  @("Greenberg Traurig", "18s", "3"),
  @("HFW", "03min 47s", "0"),
  @("BLG Law", "16s", "1")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  # Prefix column C with an apostrophe so numeric-looking registration
  # counts ("1", "2", "0", "12", ...) are stored as text, matching the
  # shared-string (t="s") cell type used throughout this column.
  $ws.Cells.Item($r, 3).Value = "'" + $row[2]
}

$lastRow = $startRow + $data.Count - 1

# Remove any leftover rows from the previous (longer) table.
$deleteRange = [string]($lastRow + 1) + ":26"
$ws.Rows($deleteRange).Delete()

# The apostrophe-prefix trick marks column C cells with a "quote prefix"
# cell style, which would otherwise diverge from the plain style used by
# columns A and B. Copy the formatting (borders/alignment/style) from
# column A onto column C to normalize the style back in line with the
# rest of the table, while keeping the text values intact.
$ws.Range("A" + $startRow).Copy() | Out-Null
$ws.Range("C" + $startRow + ":C" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
